# Auto-generated Excel COM-interop script
# Updates market-derived columns (H-N) for specific rows across multiple sheets
# reflecting a scheduled data refresh (currentAveragePrice*, LevePrice*, LeveProfit*).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3069.4
$ws.Range("I19").Value = 5615.364
$ws.Range("J19").Value = 1595.421
$ws.Range("K19").Value = 5615.364
$ws.Range("L19").Value = 1595.421
$ws.Range("M19").Value = -5440.364
$ws.Range("N19").Value = -1945.421

$ws.Range("H28").Value = 427.42856
$ws.Range("I28").Value = 427.42856
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 427.42856
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 57.57144
$ws.Range("N28").ClearContents()

$ws.Range("H64").Value = 4633.241
$ws.Range("I64").Value = 4374.9443
$ws.Range("J64").Value = 5055.909
$ws.Range("K64").Value = 4374.9443
$ws.Range("L64").Value = 5055.909
$ws.Range("M64").Value = -4126.9443
$ws.Range("N64").Value = -5551.909

$ws.Range("H67").Value = 4633.241
$ws.Range("I67").Value = 4374.9443
$ws.Range("J67").Value = 5055.909
$ws.Range("K67").Value = 4374.9443
$ws.Range("L67").Value = 5055.909
$ws.Range("M67").Value = -3516.9443
$ws.Range("N67").Value = -6771.909

$ws.Range("H88").Value = 4928.3335
$ws.Range("I88").Value = 1800.75
$ws.Range("J88").Value = 7430.4
$ws.Range("K88").Value = 1800.75
$ws.Range("L88").Value = 7430.4
$ws.Range("M88").Value = -1394.75
$ws.Range("N88").Value = -8242.4

$ws.Range("H91").Value = 4928.3335
$ws.Range("I91").Value = 1800.75
$ws.Range("J91").Value = 7430.4
$ws.Range("K91").Value = 1800.75
$ws.Range("L91").Value = 7430.4
$ws.Range("M91").Value = -396.75
$ws.Range("N91").Value = -10238.4

$ws.Range("H111").Value = 1006.3125
$ws.Range("I111").Value = 457.25
$ws.Range("J111").Value = 1189.3334
$ws.Range("K111").Value = 1371.75
$ws.Range("L111").Value = 3568.0002
$ws.Range("M111").Value = 1695.25
$ws.Range("N111").Value = -9702.0002

$ws.Range("H116").Value = 40965.75
$ws.Range("I116").Value = 58434.844
$ws.Range("J116").Value = 4086.5557
$ws.Range("K116").Value = 58434.844
$ws.Range("L116").Value = 4086.5557
$ws.Range("M116").Value = -54992.844
$ws.Range("N116").Value = -10970.5557

$ws.Range("H121").Value = 1960.5
$ws.Range("I121").Value = 597.5
$ws.Range("J121").Value = 2301.25
$ws.Range("K121").Value = 1792.5
$ws.Range("L121").Value = 6903.75
$ws.Range("M121").Value = -45.5
$ws.Range("N121").Value = -10397.75

$ws.Range("H138").Value = 1971.7872
$ws.Range("I138").Value = 2093.3125
$ws.Range("J138").Value = 1909.0646
$ws.Range("K138").Value = 6279.9375
$ws.Range("L138").Value = 5727.1938
$ws.Range("M138").Value = -1139.9375
$ws.Range("N138").Value = -16007.1938

$ws.Range("H141").Value = 10423.5
$ws.Range("I141").Value = 1231.3334
$ws.Range("J141").Value = 38000
$ws.Range("K141").Value = 3694.0002
$ws.Range("L141").Value = 114000
$ws.Range("M141").Value = 1485.9998
$ws.Range("N141").Value = -124360


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1178981.1
$ws.Range("I32").Value = 1194109.4
$ws.Range("J32").Value = 120000
$ws.Range("K32").Value = 1194109.4
$ws.Range("L32").Value = 120000
$ws.Range("M32").Value = -1193822.4
$ws.Range("N32").Value = -120574

$ws.Range("H45").Value = 1247.2307
$ws.Range("I45").Value = 678
$ws.Range("J45").Value = 2158
$ws.Range("K45").Value = 678
$ws.Range("L45").Value = 2158
$ws.Range("M45").Value = -301
$ws.Range("N45").Value = -2912

$ws.Range("H132").Value = 27532.41
$ws.Range("I132").Value = 35723.758
$ws.Range("J132").Value = 3777.5
$ws.Range("K132").Value = 107171.274
$ws.Range("L132").Value = 11332.5
$ws.Range("M132").Value = -104641.274
$ws.Range("N132").Value = -16392.5


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1051.6957
$ws.Range("I107").Value = 1165.8889
$ws.Range("J107").Value = 640.6
$ws.Range("K107").Value = 1165.8889
$ws.Range("L107").Value = 640.6
$ws.Range("M107").Value = 754.1111000000001
$ws.Range("N107").Value = -4480.6


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 440.7
$ws.Range("I5").Value = 150
$ws.Range("K5").Value = 150
$ws.Range("M5").Value = -38

$ws.Range("H16").Value = 872.17645
$ws.Range("I16").Value = 808.63635
$ws.Range("J16").Value = 988.6667
$ws.Range("K16").Value = 808.63635
$ws.Range("L16").Value = 988.6667
$ws.Range("M16").Value = -521.63635
$ws.Range("N16").Value = -1562.6667

$ws.Range("H99").Value = 113511.555
$ws.Range("I99").Value = 201960.8
$ws.Range("K99").Value = 201960.8
$ws.Range("M99").Value = -200462.8

$ws.Range("H113").Value = 872.17645
$ws.Range("I113").Value = 808.63635
$ws.Range("J113").Value = 988.6667
$ws.Range("K113").Value = 808.63635
$ws.Range("L113").Value = 988.6667
$ws.Range("M113").Value = 1361.36365
$ws.Range("N113").Value = -5328.6667

$ws.Range("H122").Value = 2763.3635
$ws.Range("I122").Value = 3249.875
$ws.Range("J122").Value = 1466
$ws.Range("K122").Value = 9749.625
$ws.Range("L122").Value = 4398
$ws.Range("M122").Value = -7299.625
$ws.Range("N122").Value = -9298

$ws.Range("H126").Value = 113511.555
$ws.Range("I126").Value = 201960.8
$ws.Range("K126").Value = 605882.3999999999
$ws.Range("M126").Value = -603412.3999999999

$ws.Range("H134").Value = 26317568
$ws.Range("I134").Value = 45455500
$ws.Range("J134").Value = 2912.5
$ws.Range("K134").Value = 136366500
$ws.Range("L134").Value = 8737.5
$ws.Range("M134").Value = -136363965
$ws.Range("N134").Value = -13807.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 10000
$ws.Range("I14").Value = 10000
$ws.Range("K14").Value = 30000
$ws.Range("M14").Value = -29827

$ws.Range("H63").Value = 1412
$ws.Range("I63").Value = 1412
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4236
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3487
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 1412
$ws.Range("I66").Value = 1412
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12708
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -8964
$ws.Range("N66").ClearContents()

$ws.Range("H68").Value = 1064.5692
$ws.Range("I68").Value = 588.925
$ws.Range("J68").Value = 1825.6
$ws.Range("K68").Value = 1766.775
$ws.Range("L68").Value = 5476.799999999999
$ws.Range("M68").Value = -955.7749999999999
$ws.Range("N68").Value = -7098.799999999999

$ws.Range("H71").Value = 1064.5692
$ws.Range("I71").Value = 588.925
$ws.Range("J71").Value = 1825.6
$ws.Range("K71").Value = 5300.325
$ws.Range("L71").Value = 16430.4
$ws.Range("M71").Value = -1244.325
$ws.Range("N71").Value = -24542.4


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1237.25
$ws.Range("I61").Value = 1164.7
$ws.Range("J61").Value = 1600
$ws.Range("K61").Value = 1164.7
$ws.Range("L61").Value = 1600
$ws.Range("M61").Value = -962.7
$ws.Range("N61").Value = -2004

$ws.Range("H113").Value = 1237.25
$ws.Range("I113").Value = 1164.7
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 1164.7
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = 1005.3
$ws.Range("N113").Value = -5940

$ws.Range("H122").Value = 474.66666
$ws.Range("I122").Value = 562
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 1686
$ws.Range("L122").Value = 900
$ws.Range("M122").Value = 764
$ws.Range("N122").Value = -5800


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 543.1818
$ws.Range("I107").Value = 502.4
$ws.Range("J107").Value = 577.1667
$ws.Range("K107").Value = 1507.2
$ws.Range("L107").Value = 1731.5001
$ws.Range("M107").Value = 412.8000000000002
$ws.Range("N107").Value = -5571.5001

$ws.Range("H126").Value = 2656.6155
$ws.Range("I126").Value = 1510.8
$ws.Range("K126").Value = 4532.4
$ws.Range("M126").Value = -2062.4

$ws.Range("H136").Value = 17001664
$ws.Range("I136").Value = 20221430
$ws.Range("J136").Value = 6939897
$ws.Range("K136").Value = 60664290
$ws.Range("L136").Value = 20819691
$ws.Range("M136").Value = -60661740
$ws.Range("N136").Value = -20824791

